$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (blank separator row) gets a red fill first, so the new "red" fill
# lands before the new "yellow" fill in the workbook's fill table - matching
# the order they appear in the target styles.xml (fillId 33 = red, 34 = yellow).
$ws.Range("A3:G3").Interior.Color = 255

# Row 2 becomes a filled-in "Example" row, shown above the real data rows.
$ws.Range("A2").Value = 501
$ws.Range("B2").Value = "May Thi Nghe"
$ws.Range("C2").Value = 132
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = "Example"

# Row 2 fill - yellow (FFFFFF00)
$ws.Range("A2:G2").Interior.Color = 65535
$ws.Range("H2").Interior.Color = 65535

# "May Thi Nghe" gets a small gray Arial label font
$ws.Range("B2").Font.Name = "Arial"
$ws.Range("B2").Font.Size = 10
$ws.Range("B2").Font.Color = 4473924

# Move the saved cursor/selection to A4, mirroring the edited workbook
$ws.Range("A4").Select()
